$wb = $excel.ActiveWorkbook

# --- Insert a new "data_collection_mode list" sheet right after "ion_mobility list" ---
$afterSheet = $wb.Worksheets.Item("ion_mobility list")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "data_collection_mode list"
$newSheet.Range("A1").Value = "DDA"
$newSheet.Range("A2").Value = "DIA"

# --- Add a matching data validation on the main sheet's data_collection_mode column (Y) ---
$main = $wb.Worksheets.Item("Export as TSV")
$range = $main.Range("Y2:Y1048576")
$range.Validation.Add(3, 1, 1, "='data_collection_mode list'!`$A`$1:`$A`$2")
$range.Validation.ErrorTitle = "Value must come from list"
$range.Validation.ErrorMessage = "Value must be one of: DDA / DIA."
$range.Validation.ShowInput = $true
$range.Validation.ShowError = $true

# --- Restore original active sheet/selection ---
$main.Activate() | Out-Null
$main.Range("A1").Select() | Out-Null
